$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateFmt = "[$-409]d/mmm/yyyy;@"
$amtFmt  = """₹""#,##0;""₹""\-#,##0"

# Row 346
$ws.Range("A346").Value = 44805
$ws.Range("A346").NumberFormat = $dateFmt
$ws.Range("B346").Value = "KA50M1697"
$ws.Range("C346").Value = "VERNA"
$ws.Range("D346").Value = "PMS"
$ws.Range("E346").Value = "WORK DONE DELIVERED"
$ws.Range("F346").Value = 19771
$ws.Range("F346").NumberFormat = $amtFmt

# Row 347
$ws.Range("A347").Value = 44805
$ws.Range("A347").NumberFormat = $dateFmt
$ws.Range("B347").Value = "KA03MS4131"
$ws.Range("C347").Value = "VERNA"
$ws.Range("D347").Value = "PMS"
$ws.Range("E347").Value = "WORK DONE DELIVERED"
$ws.Range("F347").Value = 4283
$ws.Range("F347").NumberFormat = $amtFmt
$ws.Range("G347").Value = "G PAY"

# Row 348
$ws.Range("A348").Value = 44805
$ws.Range("A348").NumberFormat = $dateFmt
$ws.Range("B348").Value = "KA01ML2754"
$ws.Range("C348").Value = "ECOSPORT"
$ws.Range("D348").Value = "RUNNING REPAIR"
$ws.Range("E348").Value = "WORK DONE DELIVERED"
$ws.Range("F348").Value = 1850
$ws.Range("F348").NumberFormat = $amtFmt
$ws.Range("G348").Value = "CREDIT"

# Row 349
$ws.Range("A349").Value = 44774
$ws.Range("A349").NumberFormat = $dateFmt
$ws.Range("B349").Value = "KA03MQ6297"
$ws.Range("C349").Value = "FIGO"
$ws.Range("D349").Value = "PMS"
$ws.Range("E349").Value = "WORK DONE DELIVERED"
$ws.Range("F349").Value = 4322
$ws.Range("F349").NumberFormat = $amtFmt
$ws.Range("G349").Value = "P PAY"

# Row 350
$ws.Range("A350").Value = 44806
$ws.Range("A350").NumberFormat = $dateFmt
$ws.Range("B350").Value = "KA02MA7199"
$ws.Range("C350").Value = "SCORPIO"
$ws.Range("D350").Value = "RUNNING REPAIR"
$ws.Range("E350").Value = "WORK DONE DELIVERED"
$ws.Range("F350").Value = 400
$ws.Range("F350").NumberFormat = $amtFmt
$ws.Range("G350").Value = "P PAY"

# Row 351
$ws.Range("A351").Value = 44806
$ws.Range("A351").NumberFormat = $dateFmt
$ws.Range("B351").Value = "KA03ML5436"
$ws.Range("C351").Value = "H-CITY"
$ws.Range("D351").Value = "GENERAL CHECKUP         WW"
$ws.Range("E351").Value = "WORK IN PROGRESS"

# Best-effort: restore the scrolled viewport position (topLeftCell) and
# final selection, matching the author's on-screen state after the edit.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 334
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E346").Select() | Out-Null
